# Replace old localhost URLs with the production domain across the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$oldPrefix = "http://localhost:3000"
$newPrefix = "https://t-h-logistics.com"

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count
$colCount = $usedRange.Columns.Count

for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value()
        if ($val -ne $null -and $val -is [string] -and $val.StartsWith($oldPrefix)) {
            $cell.Value = $newPrefix + $val.Substring($oldPrefix.Length)
        }
    }
}
